$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B8").Value = "In Translation"
$zhcn.Range("B9").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B8").Value = "In Translation"
$dede.Range("B9").Value = "In Translation"
